$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A23 holds a date-like literal string ("05/08/2025"). A plain .Value
# assignment gets auto-parsed into a date serial by Excel's text-to-number
# inference, so force text entry (NumberFormat "@") then drop back to the
# default "Normal" style once the literal text is safely stored.
$dateCell = $ws.Cells.Item(23, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "05/08/2025"
$dateCell.Style = "Normal"

$ws.Range("B23").Value = "Melgar"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = "Juan Pablo II"
$ws.Range("F23").Value = "D"
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 1.38
$ws.Range("L23").Value = 0.17
$ws.Range("M23").Value = 21
$ws.Range("N23").Value = 5
$ws.Range("O23").Value = 5
$ws.Range("P23").Value = 2
